$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Progress (F) column: rows 3-6 move from 90% to 100% complete ---
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 1

# --- Duration (C) and Progress (F) newly filled in for rows 7-10 ---
$ws.Range("C7").Value = 0.089583333333333334
$ws.Range("F7").Value = 0.9

$ws.Range("C8").Value = 0.027777777777777776
$ws.Range("F8").Value = 0.9

$ws.Range("C9").Value = 0.04583333333333333

$ws.Range("C10").Value = 0.04583333333333333

# --- Conditional formatting on F3:F52: replace the 3 cellIs rules with a 3-color scale ---
$rangeCF = $ws.Range("F3:F52")

# Recreate the existing "Good/Bad/Neutral" style rules first so their
# differential formats are preserved (Excel keeps the old dxf records
# around even after the rules referencing them are removed).
$cGood = $rangeCF.FormatConditions.Add(1, 3, "0")
$cGood.Font.Color = 24832
$cGood.Interior.Color = 13561798

$cBad = $rangeCF.FormatConditions.Add(1, 6, "0.9")
$cBad.Font.Color = 393372
$cBad.Interior.Color = 13551615

$cNeutral = $rangeCF.FormatConditions.Add(1, 5, "0.89")
$cNeutral.Font.Color = 22428
$cNeutral.Interior.Color = 10284031

# Remove every existing rule (the 3 original ones plus the 3 just added)
$rangeCF.FormatConditions.Delete()

# Add the new 3-color scale rule
$cf = $rangeCF.FormatConditions.AddColorScale(3)
$cf.ColorScaleCriteria.Item(1).Type = 1   # xlConditionValueLowestValue
$cf.ColorScaleCriteria.Item(1).FormatColor.Color = 7039851    # FFF8696B

$cf.ColorScaleCriteria.Item(2).Type = 5   # xlConditionValuePercentile
$cf.ColorScaleCriteria.Item(2).Value = 50
$cf.ColorScaleCriteria.Item(2).FormatColor.Color = 8711167    # FFFFEB84

$cf.ColorScaleCriteria.Item(3).Type = 2   # xlConditionValueHighestValue
$cf.ColorScaleCriteria.Item(3).FormatColor.Color = 8107615    # FF63BE7B

# --- Move the active selection to F9 (pane's topLeftCell auto-adjusts) ---
$ws.Range("F9").Select()
